$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 02_HW_Components: append "to buy" section
# ----------------------------------------------------------------------
$wsHw = $wb.Worksheets.Item(2)
$wsHw.Range("A16").Value = "to buy"
$wsHw.Range("B16").Value = "small smd shotky diodes  (100mA)"
$wsHw.Range("B17").Value = "supercap"
$wsHw.Range("J6").Select()

# ----------------------------------------------------------------------
# Rename Sheet3 -> "Power source", fill in voltage-regulator calculations
# ----------------------------------------------------------------------
$wsPower = $wb.Worksheets.Item(3)
$wsPower.Name = "Power source"

$wsPower.Range("A3").Value = "Calculations for XC9264B75DER-G:"

$wsPower.Range("B5").Value = "V.in"
$wsPower.Range("C5").Value = 6
$wsPower.Range("D5").Value = "V"
$wsPower.Range("E5").Value = "Cfb"
$wsPower.Range("F5").Formula = "=1/(2*PI() * C8 * C9* 1000) * 1000000000"
$wsPower.Range("F5").Font.Bold = $true
$wsPower.Range("F5").Font.Size = 11
$wsPower.Range("F5").Font.ThemeColor = 5
$wsPower.Range("G5").Value = "nF"

$wsPower.Range("B6").Value = "V.out"
$wsPower.Range("C6").Value = 3.3
$wsPower.Range("D6").Value = "V"
$wsPower.Range("E6").Value = "Rfb2"
$wsPower.Range("F6").Formula = "=(0.75*C9)/(C6-0.75)"
$wsPower.Range("F6").Font.Bold = $true
$wsPower.Range("F6").Font.Size = 11
$wsPower.Range("F6").Font.ThemeColor = 5
$wsPower.Range("G6").Value = "kOhm"

$wsPower.Range("B7").Value = "L"
$wsPower.Range("C7").Value = 2.2000000000000002
$wsPower.Range("D7").Value = "uH"

$wsPower.Range("B8").Value = "fzfb"
$wsPower.Range("C8").Value = 5000
$wsPower.Range("D8").Value = "Hz"

$wsPower.Range("B9").Value = "Rfb1"
$wsPower.Range("C9").Value = 100
$wsPower.Range("D9").Value = "kOhm"

$wsPower.Range("B13").Value = "Vfb"
$wsPower.Range("C13").Value = 0.75
$wsPower.Range("D13").Value = 0.73899999999999999
$wsPower.Range("E13").Value = 0.76100000000000001

$wsPower.Range("B14").Value = "Rfb1"
$wsPower.Range("C14").Value = 100

$wsPower.Range("B15").Value = "Rfb2"
$wsPower.Range("C15").Value = 33

$wsPower.Range("B16").Value = "Vout"
$wsPower.Range("C16").Formula = "=C13*(`$C`$14+`$C`$15)/`$C`$15"
$wsPower.Range("C16").Font.Bold = $true
$wsPower.Range("C16").Font.Size = 11
$wsPower.Range("C16").Font.ThemeColor = 3
$wsPower.Range("D16").Formula = "=D13*(`$C`$14+`$C`$15)/`$C`$15"
$wsPower.Range("D16").Font.Bold = $true
$wsPower.Range("D16").Font.Size = 11
$wsPower.Range("D16").Font.ThemeColor = 3
$wsPower.Range("E16").Formula = "=E13*(`$C`$14+`$C`$15)/`$C`$15"
$wsPower.Range("E16").Font.Bold = $true
$wsPower.Range("E16").Font.Size = 11
$wsPower.Range("E16").Font.ThemeColor = 3

$wsPower.Range("F5").Select()

# ----------------------------------------------------------------------
# New sheet "ExtFlash" with datasheet timing/power parameters
# (values entered in the same order the original author typed them, so
# the shared-string table comes out in the same sequence)
# ----------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsFlash = $wb.Worksheets.Add($null, $lastSheet)
$wsFlash.Name = "ExtFlash"

$wsFlash.Range("B5").Value = "write power req. "
$wsFlash.Range("D5").Value = "50mA"

$wsFlash.Range("B2").Value = "Von min"
$wsFlash.Range("D2").Value = "2.7V"

$wsFlash.Range("B3").Value = "Vcut"
$wsFlash.Range("D3").Value = "2.4V"

$wsFlash.Range("B6").Value = "read power req."
$wsFlash.Range("D6").Value = "15mA"
$wsFlash.Range("E6").Value = "with 20MHz clock"

$wsFlash.Range("B7").Value = "stdby"
$wsFlash.Range("D7").Value = "30uA"

$wsFlash.Range("B9").Value = "T reset pulse"
$wsFlash.Range("D9").Value = "200ns"
$wsFlash.Range("E9").Value = "min"

$wsFlash.Range("E10").Value = "min"
$wsFlash.Range("B10").Value = "T reset procedure"
$wsFlash.Range("D10").Value = "100us"

$wsFlash.Range("B12").Value = "T page prog 256"
$wsFlash.Range("D12").Value = "450us "
$wsFlash.Range("E12").Value = "1350us"

$wsFlash.Range("E11").Value = "max:"
$wsFlash.Range("D11").Value = "typ:"

$wsFlash.Range("B13").Value = "T byte prog 1st"
$wsFlash.Range("D13").Value = "75us"
$wsFlash.Range("E13").Value = "90us"

$wsFlash.Range("B14").Value = "T byte prog next"
$wsFlash.Range("D14").Value = "10us"
$wsFlash.Range("E14").Value = "30us"

$wsFlash.Range("B15").Value = "T err sector 4k"
$wsFlash.Range("D15").Value = "65ms"
$wsFlash.Range("E15").Value = "320ms"

$wsFlash.Range("B16").Value = "T err 1/2bk 32kB"
$wsFlash.Range("D16").Value = "300ms"
$wsFlash.Range("E16").Value = "600ms"

$wsFlash.Range("B17").Value = "T err bk 64kB"
$wsFlash.Range("D17").Value = "450ms"
$wsFlash.Range("E17").Value = "1150ms"

$wsFlash.Range("B18").Value = "T err all"
$wsFlash.Range("D18").Value = "55s"
$wsFlash.Range("E18").Value = "150s"

$wsFlash.Range("B19").Select()

# ExtFlash is the sheet that is active/selected when the file is saved
$wsFlash.Activate()
